$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A17:E17").ClearContents()
$ws.Range("A17:E17").ClearFormats()
